# Refresh ligand/receptor TPM-derived metrics (Inha -> Acvr2a) with updated source values.
# Mirrors the NATMI recompute: detection counts/rates, average & total expression,
# derived specificities, and edge weights for rows 2-9 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.6944570264822121
$ws.Range("J2").Value = 0.6944570264822121
$ws.Range("M2").Value = 16.28844733333333
$ws.Range("N2").Value = 48.865342
$ws.Range("O2").Value = 0.2176904746803693
$ws.Range("P2").Value = 0.2176904746803693
$ws.Range("Q2").Value = 4.361959324147556
$ws.Range("R2").Value = 39.257633917328
$ws.Range("S2").Value = 0.1511766797400306
$ws.Range("T2").Value = 0.1511766797400306

# Row 3
$ws.Range("I3").Value = 0.6944570264822121
$ws.Range("J3").Value = 0.6944570264822121
$ws.Range("M3").Value = 27.61090666666666
$ws.Range("N3").Value = 82.83272
$ws.Range("O3").Value = 0.3690119294748028
$ws.Range("P3").Value = 0.3690119294748029
$ws.Range("Q3").Value = 7.394053547164444
$ws.Range("R3").Value = 66.54648192447999
$ws.Range("S3").Value = 0.2562629272795353
$ws.Range("T3").Value = 0.2562629272795354

# Row 4
$ws.Range("I4").Value = 0.6944570264822121
$ws.Range("J4").Value = 0.6944570264822121
$ws.Range("M4").Value = 26.266325
$ws.Range("N4").Value = 78.798975
$ws.Range("O4").Value = 0.3510419771967738
$ws.Range("P4").Value = 0.3510419771967739
$ws.Range("Q4").Value = 7.033981747933333
$ws.Range("R4").Value = 63.30583573139999
$ws.Range("S4").Value = 0.243783567654508
$ws.Range("T4").Value = 0.2437835676545081

# Row 5
$ws.Range("I5").Value = 0.6944570264822121
$ws.Range("J5").Value = 0.6944570264822121
$ws.Range("M5").Value = 4.658207333333333
$ws.Range("N5").Value = 13.974622
$ws.Range("O5").Value = 0.06225561864805391
$ws.Range("P5").Value = 0.06225561864805392
$ws.Range("Q5").Value = 1.247443080094222
$ws.Range("R5").Value = 11.226987720848
$ws.Range("S5").Value = 0.04323385180813807
$ws.Range("T5").Value = 0.04323385180813808

# Row 6
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.1178226666666667
$ws.Range("H6").Value = 0.353468
$ws.Range("I6").Value = 0.3055429735177879
$ws.Range("J6").Value = 0.3055429735177879
$ws.Range("M6").Value = 16.28844733333333
$ws.Range("N6").Value = 48.865342
$ws.Range("O6").Value = 0.2176904746803693
$ws.Range("P6").Value = 0.2176904746803693
$ws.Range("Q6").Value = 1.919148300672889
$ws.Range("R6").Value = 17.272334706056
$ws.Range("S6").Value = 0.06651379494033877
$ws.Range("T6").Value = 0.06651379494033877

# Row 7
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.1178226666666667
$ws.Range("H7").Value = 0.353468
$ws.Range("I7").Value = 0.3055429735177879
$ws.Range("J7").Value = 0.3055429735177879
$ws.Range("M7").Value = 27.61090666666666
$ws.Range("N7").Value = 82.83272
$ws.Range("O7").Value = 0.3690119294748028
$ws.Range("P7").Value = 0.3690119294748029
$ws.Range("Q7").Value = 3.253190652551111
$ws.Range("R7").Value = 29.27871587296
$ws.Range("S7").Value = 0.1127490021952675
$ws.Range("T7").Value = 0.1127490021952675

# Row 8
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.1178226666666667
$ws.Range("H8").Value = 0.353468
$ws.Range("I8").Value = 0.3055429735177879
$ws.Range("J8").Value = 0.3055429735177879
$ws.Range("M8").Value = 26.266325
$ws.Range("N8").Value = 78.798975
$ws.Range("O8").Value = 0.3510419771967738
$ws.Range("P8").Value = 0.3510419771967739
$ws.Range("Q8").Value = 3.094768455033333
$ws.Range("R8").Value = 27.8529160953
$ws.Range("S8").Value = 0.1072584095422658
$ws.Range("T8").Value = 0.1072584095422658

# Row 9
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.1178226666666667
$ws.Range("H9").Value = 0.353468
$ws.Range("I9").Value = 0.3055429735177879
$ws.Range("J9").Value = 0.3055429735177879
$ws.Range("M9").Value = 4.658207333333333
$ws.Range("N9").Value = 13.974622
$ws.Range("O9").Value = 0.06225561864805391
$ws.Range("P9").Value = 0.06225561864805392
$ws.Range("Q9").Value = 0.5488424098995556
$ws.Range("R9").Value = 4.939581689096
$ws.Range("S9").Value = 0.01902176683991584
$ws.Range("T9").Value = 0.01902176683991584
